# Updated symbol list (price / volume refresh + BOLO<->CoinbaseStockToken
# row swap at rows 48-49). Price/volume cells hold numeric-looking text
# (e.g. "307.61", "-0.33%") in the source data, so each is assigned with a
# leading apostrophe to force Excel to keep it as text instead of
# auto-converting it to a Number/Percentage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.61"
$ws.Range("E2").Value = "'-0.33%"
$ws.Range("D3").Value = "'41.00"
$ws.Range("E3").Value = "'0.02%"
$ws.Range("D4").Value = "'5.044"
$ws.Range("E4").Value = "'-1.04%"
$ws.Range("D5").Value = "'0.07609"
$ws.Range("E5").Value = "'-2.94%"
$ws.Range("D6").Value = "'1.601"
$ws.Range("E6").Value = "'-5.25%"
$ws.Range("D7").Value = "'0.9063"
$ws.Range("E7").Value = "'-1.49%"
$ws.Range("E8").Value = "'-5.10%"
$ws.Range("D9").Value = "'0.1021"
$ws.Range("E9").Value = "'-4.28%"
$ws.Range("D10").Value = "'0.1763"
$ws.Range("E10").Value = "'0.17%"
$ws.Range("D11").Value = "'0.09081"
$ws.Range("E11").Value = "'0.43%"
$ws.Range("D12").Value = "'0.04349"
$ws.Range("E12").Value = "'-1.42%"
$ws.Range("D13").Value = "'0.1052"
$ws.Range("E13").Value = "'-0.48%"
$ws.Range("D14").Value = "'0.001256"
$ws.Range("E14").Value = "'-3.06%"
$ws.Range("D15").Value = "'0.005880"
$ws.Range("E15").Value = "'0.88%"
$ws.Range("D16").Value = "'3.354"
$ws.Range("E16").Value = "'-0.74%"
$ws.Range("D17").Value = "'4.270"
$ws.Range("E17").Value = "'-1.07%"
$ws.Range("D18").Value = "'0.3269"
$ws.Range("E18").Value = "'-2.95%"
$ws.Range("D19").Value = "'6.814"
$ws.Range("E19").Value = "'-5.31%"
$ws.Range("D20").Value = "'0.1358"
$ws.Range("E20").Value = "'-2.09%"
$ws.Range("D21").Value = "'0.2728"
$ws.Range("E21").Value = "'-2.78%"
$ws.Range("D22").Value = "'0.04182"
$ws.Range("E22").Value = "'0.25%"
$ws.Range("D23").Value = "'0.001228"
$ws.Range("E23").Value = "'0.63%"
$ws.Range("D24").Value = "'0.004064"
$ws.Range("E24").Value = "'-1.75%"
$ws.Range("D25").Value = "'0.0001303"
$ws.Range("E25").Value = "'6.22%"
$ws.Range("D26").Value = "'0.0003013"
$ws.Range("E26").Value = "'0.56%"
$ws.Range("D38").Value = "'0.02397"
$ws.Range("E38").Value = "'-1.24%"
$ws.Range("D39").Value = "'0.05181"
$ws.Range("E39").Value = "'-1.33%"
$ws.Range("D40").Value = "'0.007796"
$ws.Range("E40").Value = "'-2.05%"
$ws.Range("D41").Value = "'0.1303"
$ws.Range("E41").Value = "'-3.55%"
$ws.Range("D42").Value = "'0.007089"
$ws.Range("E42").Value = "'-5.33%"
$ws.Range("D43").Value = "'0.001922"
$ws.Range("E43").Value = "'-4.95%"
$ws.Range("D44").Value = "'0.007492"
$ws.Range("E44").Value = "'-7.55%"
$ws.Range("D45").Value = "'0.3349"
$ws.Range("E45").Value = "'-0.56%"
$ws.Range("D46").Value = "'0.00006358"
$ws.Range("E46").Value = "'-6.08%"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("E47").Value = "'-0.37%"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "'0.006574"
$ws.Range("E48").Value = "'105.51%"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").Value = "'0.004409"
$ws.Range("E49").Value = "'6.90%"
$ws.Range("D50").Value = "'0.00002104"
$ws.Range("E50").Value = "'-0.37%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'-0.37%"
